$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: replace Sung Reichel with Agatha Bauch
$ws.Range("B7").Value = "agathabauch@bsgdulpk.mailosaur.net"
$ws.Range("D7").Value = "Agatha"
$ws.Range("E7").Value = "Bauch"

# Row 8: replace Miles Hodkiewicz with Tatiana Wehner
$ws.Range("B8").Value = "tatianawehner@bsgdulpk.mailosaur.net"
$ws.Range("D8").Value = "Tatiana"
$ws.Range("E8").Value = "Wehner"
